$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at the top of the data (row 2), pushing every existing
# price-history row (and its hyperlink) down by one. This is how the sheet
# tracks a new day: the whole history shifts down and a new "latest" row is
# added at the top.
$ws.Rows.Item(2).Insert()

# Row 3 now holds what used to be row 2 (the previous "latest" entry, dated
# 20-01-2026 with a basic price of 338.5 from the 20-01-2026 circular). The
# brand-new row 2 repeats that same price/circular info, just one day later.
$ws.Range("B2:F2").Value2 = $ws.Range("B3:F3").Value2
$ws.Range("A2").Value2 = "21-01-2026"

# Match the new row's formatting (borders/alignment/number format) to the
# rest of the data rows.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
